$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - account holder first name
$ws.Range("C2").Value = "Hartmut"

# Row 3 - card number (purely numeric-looking text, must stay text) + account holder last name
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "2570314725427075"
$ws.Range("Z1").Copy()
$ws.Range("B3").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("C3").Value = "Mohaupt"

# Row 5 - opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 02.07.2024"

# Row 6 - transaction 1
$ws.Range("B6").Value = "03.07."
$ws.Range("C6").Value = "04.07."
$ws.Range("D6").Value = "KARTENZ./03.07 REWE RO"
$ws.Range("E6").Value = "92,68-"

# Row 7 - transaction 2
$ws.Range("B7").Value = "04.07."
$ws.Range("C7").Value = "05.07."
$ws.Range("D7").Value = "KARTENZAHLUNG ARAL TANKSTELLE"
$ws.Range("E7").Value = "89,84-"

# Row 8 - transaction 3
$ws.Range("B8").Value = "07.07."
$ws.Range("C8").Value = "08.07."
$ws.Range("D8").Value = "RECHNUNG VODAFONE GMBH 66926558"
$ws.Range("E8").Value = "38,44-"

# Row 9 - transaction removed entirely, row becomes blank (E9:F9 is a merged cell)
$ws.Range("B9:D9").ClearContents()
$ws.Range("E9:F9").ClearContents()
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").VerticalAlignment = -4108
$ws.Range("E9").WrapText = $true

# Row 12 - closing balance
$ws.Range("D12").Value = "KONTOSTAND AM 12.07.2024"
$ws.Range("E12").Value = "220,96-"

# Row 13 - next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 18.07.2024"
